# Weekly fruit/vegetable price refresh: the per-record fields (date, volume,
# price range, commercialization unit, origin, $/Kg and Kg-or-unit basis) are
# reshuffled across the existing rows of the "Camote" (sweet potato) table for
# Vega Modelo de Temuco. The market/region/category/quality columns (A, B, C,
# E, F, G, H, I, R) are identical for every record in this sheet and stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = 44424; J = 30; K = 20000; L = 20000; M = 20000; N = "`$/caja 15 kilos granel"; O = "Región de Arica y Parinacota"; P = 1333; Q = 15 },
    @{ Row = 3; D = 44294; J = 5; K = 20000; L = 20000; M = 20000; N = "`$/caja 15 kilos granel"; O = "Perú"; P = 1333; Q = 15 },
    @{ Row = 4; D = 44511; J = 50; K = 20000; L = 20000; M = 20000; N = "`$/malla 20 kilos"; O = "Perú"; P = 1000; Q = 20 },
    @{ Row = 5; D = 44364; J = 15; K = 20000; L = 20000; M = 20000; N = "`$/caja 15 kilos granel"; O = "Perú"; P = 1333; Q = 15 },
    @{ Row = 6; D = 44316; J = 20; K = 20000; L = 20000; M = 20000; N = "`$/caja 15 kilos granel"; O = "Región de Arica y Parinacota"; P = 1333; Q = 15 },
    @{ Row = 7; D = 44466; J = 20; K = 25000; L = 25000; M = 25000; N = "`$/caja 15 kilos granel"; O = "Perú"; P = 1667; Q = 15 },
    @{ Row = 8; D = 44385; J = 18; K = 20000; L = 20000; M = 20000; N = "`$/malla 20 kilos"; O = "Región de Arica y Parinacota"; P = 1000; Q = 20 },
    @{ Row = 9; D = 44179; J = 20; K = 20000; L = 20000; M = 20000; N = "`$/caja 15 kilos granel"; O = "Región de Arica y Parinacota"; P = 1333; Q = 15 },
    @{ Row = 10; D = 44455; J = 30; K = 20000; L = 20000; M = 20000; N = "`$/malla 20 kilos"; O = "Perú"; P = 1000; Q = 20 },
    @{ Row = 11; D = 44369; J = 20; K = 20000; L = 20000; M = 20000; N = "`$/caja 15 kilos granel"; O = "Región de Arica y Parinacota"; P = 1333; Q = 15 },
    @{ Row = 12; D = 44369; J = 20; K = 20000; L = 20000; M = 20000; N = "`$/malla 20 kilos"; O = "Región de Arica y Parinacota"; P = 1000; Q = 20 },
    @{ Row = 13; D = 44188; J = 20; K = 20000; L = 20000; M = 20000; N = "`$/caja 15 kilos granel"; O = "Región de Arica y Parinacota"; P = 1333; Q = 15 },
    @{ Row = 14; D = 44512; J = 30; K = 20000; L = 20000; M = 20000; N = "`$/malla 20 kilos"; O = "Perú"; P = 1000; Q = 20 },
    @{ Row = 15; D = 44186; J = 20; K = 20000; L = 20000; M = 20000; N = "`$/caja 15 kilos granel"; O = "Región de Arica y Parinacota"; P = 1333; Q = 15 },
    @{ Row = 16; D = 44438; J = 40; K = 20000; L = 20000; M = 20000; N = "`$/caja 15 kilos granel"; O = "Región de Arica y Parinacota"; P = 1333; Q = 15 },
    @{ Row = 17; D = 44497; J = 30; K = 20000; L = 20000; M = 20000; N = "`$/caja 15 kilos granel"; O = "Perú"; P = 1333; Q = 15 },
    @{ Row = 18; D = 44497; J = 40; K = 20000; L = 20000; M = 20000; N = "`$/malla 20 kilos"; O = "Perú"; P = 1000; Q = 20 },
    @{ Row = 19; D = 44175; J = 20; K = 20000; L = 20000; M = 20000; N = "`$/caja 15 kilos granel"; O = "Región de Arica y Parinacota"; P = 1333; Q = 15 },
    @{ Row = 20; D = 44496; J = 30; K = 20000; L = 20000; M = 20000; N = "`$/malla 20 kilos"; O = "Perú"; P = 1000; Q = 20 },
    @{ Row = 21; D = 44498; J = 20; K = 20000; L = 20000; M = 20000; N = "`$/malla 20 kilos"; O = "Región de Arica y Parinacota"; P = 1000; Q = 20 },
    @{ Row = 22; D = 44321; J = 15; K = 25000; L = 25000; M = 25000; N = "`$/caja 15 kilos granel"; O = "Perú"; P = 1667; Q = 15 },
    @{ Row = 23; D = 44452; J = 50; K = 20000; L = 20000; M = 20000; N = "`$/malla 20 kilos"; O = "Perú"; P = 1000; Q = 20 },
    @{ Row = 24; D = 44448; J = 45; K = 20000; L = 20000; M = 20000; N = "`$/malla 20 kilos"; O = "Perú"; P = 1000; Q = 20 },
    @{ Row = 25; D = 44441; J = 40; K = 20000; L = 20000; M = 20000; N = "`$/malla 20 kilos"; O = "Perú"; P = 1000; Q = 20 },
    @{ Row = 26; D = 44341; J = 40; K = 17000; L = 18000; M = 17500; N = "`$/malla 20 kilos"; O = "Perú"; P = 875; Q = 20 },
    @{ Row = 27; D = 44329; J = 40; K = 20000; L = 20000; M = 20000; N = "`$/caja 15 kilos granel"; O = "Perú"; P = 1333; Q = 15 },
    @{ Row = 28; D = 44425; J = 10; K = 20000; L = 20000; M = 20000; N = "`$/caja 15 kilos granel"; O = "Región de Arica y Parinacota"; P = 1333; Q = 15 },
    @{ Row = 29; D = 44315; J = 30; K = 20000; L = 20000; M = 20000; N = "`$/caja 15 kilos granel"; O = "Región de Arica y Parinacota"; P = 1333; Q = 15 },
    @{ Row = 30; D = 44315; J = 30; K = 20000; L = 20000; M = 20000; N = "`$/malla 20 kilos"; O = "Región de Arica y Parinacota"; P = 1000; Q = 20 },
    @{ Row = 31; D = 44161; J = 20; K = 20000; L = 20000; M = 20000; N = "`$/caja 15 kilos granel"; O = "Región de Arica y Parinacota"; P = 1333; Q = 15 },
    @{ Row = 32; D = 44389; J = 45; K = 20000; L = 20000; M = 20000; N = "`$/caja 15 kilos granel"; O = "Región de Arica y Parinacota"; P = 1333; Q = 15 },
    @{ Row = 33; D = 44340; J = 40; K = 18000; L = 18000; M = 18000; N = "`$/malla 20 kilos"; O = "Perú"; P = 900; Q = 20 },
    @{ Row = 34; D = 44508; J = 40; K = 20000; L = 20000; M = 20000; N = "`$/caja 15 kilos granel"; O = "Perú"; P = 1333; Q = 15 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 4).Value  = $u.D   # Fecha
    $ws.Cells.Item($r, 10).Value = $u.J   # Volumen
    $ws.Cells.Item($r, 11).Value = $u.K   # Precio mínimo
    $ws.Cells.Item($r, 12).Value = $u.L   # Precio máximo
    $ws.Cells.Item($r, 13).Value = $u.M   # Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $u.N   # Unidad de comercialización
    $ws.Cells.Item($r, 15).Value = $u.O   # Origen
    $ws.Cells.Item($r, 16).Value = $u.P   # Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $u.Q   # Kg o Unidades
}
